$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
# Existing headers (A1:D1) keep referencing the original shared strings.
$ws.Range("A1").Value = "Model Name"
$ws.Range("B1").Value = "Training Accuracy"
$ws.Range("D1").Value = "Validation Accuracy"
$ws.Range("I1").Value = "Testing Accuracy "

# New "Loss" headers, added next (matches authoring order of shared strings).
$ws.Range("C1").Value = "Training Loss"
$ws.Range("E1").Value = "Validation Loss"
$ws.Range("J1").Value = "Testing Loss"

# Model name labels (column A, rows 2-4).
$ws.Range("A2").Value = "ResNet50"
$ws.Range("A3").Value = "VGG16"
$ws.Range("A4").Value = "InceptionV3"

# "Testing Time" group headers (K, M, L).
$ws.Range("K1").Value = "Testing Time"
$ws.Range("M1").Value = "Average Time per Classification"
$ws.Range("L1").Value = "Number of images"

# "Train and Validation Time" group headers (F, G, H) -- added last.
$ws.Range("F1").Value = "Train and Validation Time (s)"
$ws.Range("G1").Value = "Number of Images"
$ws.Range("H1").Value = "Time per image (tbc)"

# ---- Row 2: ResNet50 data ----
$ws.Range("B2").Value = 0.98429999999999995
$ws.Range("C2").Value = 0.0504
$ws.Range("D2").Value = 0.95
$ws.Range("E2").Value = 0.2464
$ws.Range("I2").Value = 0.9333
$ws.Range("J2").Value = 0.037

# ---- Row 3: VGG16 data ----
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.013
$ws.Range("D3").Value = 0.975
$ws.Range("E3").Value = 0.084
$ws.Range("F3").Value = 2343
$ws.Range("I3").Value = 0.9778
$ws.Range("J3").Value = 0.0007
$ws.Range("K3").Value = 113

# ---- Row 4: InceptionV3 data ----
$ws.Range("B4").Value = 0.99409999999999998
$ws.Range("C4").Value = 0.0878
$ws.Range("D4").Value = 0.9625
$ws.Range("E4").Value = 0.8078
$ws.Range("F4").Value = 474
$ws.Range("I4").Value = 0.9556
$ws.Range("J4").Value = 0.7297
$ws.Range("K4").Value = 24

# ---- Column L: constant 180, filled for rows 2-7 ----
$ws.Range("L2").Formula = "=36*5"
$ws.Range("L3:L7").Value = 180

# ---- Column M: K/L ratio. M2 stands alone; M3:M7 form a shared formula group ----
$ws.Range("M2").Formula = "=K2/L2"
$ws.Range("M3:M7").Formula = "=K3/L3"

# ---- Column widths (best effort match to authored widths; this runtime
#      quantizes ColumnWidth to 1/6-character steps, so these inputs are
#      chosen to land on the output step nearest the authored width) ----
$ws.Columns.Item(2).ColumnWidth = 15
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 17
$ws.Columns.Item(5).ColumnWidth = 17
$ws.Columns.Item(6).ColumnWidth = 24.6
$ws.Columns.Item(7).ColumnWidth = 24.6
$ws.Columns.Item(8).ColumnWidth = 24.6
$ws.Columns.Item(9).ColumnWidth = 14.8
$ws.Columns.Item(12).ColumnWidth = 15.5

# ---- Sheet view selection ----
$ws.Range("D6").Select() | Out-Null
